# Assignment1.xlsx edit: remove the "Assignment" worksheet (solver scratch sheet)
# and refresh the PO List statistics that depended on it.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Delete the "Assignment" sheet -----------------------------------------
[void]$wb.Worksheets.Item("Assignment").Delete()

# --- Refresh "PO List" stats (10/30-storey counters & dates) ---------------
$ws = $wb.Worksheets.Item("PO List")

$ws.Range("N3").Value = 13
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 44817

$ws.Range("N4").Value = 18
$ws.Range("R4").Value = 10

$ws.Range("N5").Value = 21
$ws.Range("O5").Value = 5
$ws.Range("P5").Value = 5
$ws.Range("R5").Value = 6

$ws.Range("M6").Value = 44385
$ws.Range("N6").Value = 22
$ws.Range("O6").Value = 3
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 44774
$ws.Range("R6").Value = 16

$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 3
$ws.Range("P7").Value = 3
$ws.Range("R7").Value = 21

$ws.Range("N8").Value = 20
$ws.Range("R8").Value = 3

$ws.Range("M9").Value = 44749
$ws.Range("N9").Value = 16
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 44900
$ws.Range("R9").Value = 5

$ws.Range("N10").Value = 12
$ws.Range("O10").Value = 3
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 44834
$ws.Range("R10").Value = 12

$ws.Range("M11").Value = 44869
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 3
$ws.Range("P11").Value = 3
$ws.Range("Q11").Value = 44951
$ws.Range("R11").Value = 1

$ws.Range("N12").Value = 8
$ws.Range("O12").Value = 4
$ws.Range("P12").Value = 4
$ws.Range("Q12").Value = 44433
$ws.Range("R12").Value = 22

$ws.Range("N13").Value = 7
$ws.Range("R13").Value = 2

$ws.Range("N14").Value = 19

$ws.Range("M15").Value = 44182
$ws.Range("N15").Value = 25
$ws.Range("O15").Value = 3
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 44819
$ws.Range("R15").Value = 13

$ws.Range("N16").Value = 14
$ws.Range("R16").Value = 17

$ws.Range("M17").Value = 44260
$ws.Range("N17").Value = 24
$ws.Range("O17").Value = 4
$ws.Range("P17").Value = 4
$ws.Range("Q17").Value = 44809
$ws.Range("R17").Value = 15

$ws.Range("M18").Value = 44328
$ws.Range("N18").Value = 23
$ws.Range("O18").Value = 4
$ws.Range("P18").Value = 4
$ws.Range("R18").Value = 8

$ws.Range("N19").Value = 5

$ws.Range("N20").Value = 2
$ws.Range("R20").Value = 18

$ws.Range("N21").Value = 17
$ws.Range("R21").Value = 20

$ws.Range("N22").Value = 26

$ws.Range("M23").Value = 31140
$ws.Range("O23").Value = 2
$ws.Range("P23").Value = 2
$ws.Range("R23").Value = 23

$ws.Range("N24").Value = 9
$ws.Range("O24").Value = 4
$ws.Range("P24").Value = 4
$ws.Range("R24").Value = 19

$ws.Range("N25").Value = 6
$ws.Range("R25").Value = 4

$ws.Range("N26").Value = 15
$ws.Range("R26").Value = 11

$ws.Range("N27").Value = 4
$ws.Range("R27").Value = 9

$ws.Range("N28").Value = 1
$ws.Range("Z28").Value = 1

$ws.Range("N29").Value = 11
$ws.Range("R29").Value = 6

# --- Restore view state: "PO List" active, cursor at R32 -------------------
[void]$ws.Activate()
[void]$ws.Range("R32").Select()
